$d = $word.ActiveDocument

# 1. Update the cover-letter date ("30 June 2021" -> "1 July 2021").
$d.Content.Find.Execute("30 June 2021", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1 July 2021", 2)

# 2. Italicize the species name "Coregonus artedi" in the submission
#    paragraph. This splits the original single run into three runs:
#    plain text, the italicized species name, and the remaining plain text.
$r = $d.Content
$r.Find.Execute("Coregonus artedi", $true, $false, $false, $false, $false,
                $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Font.Italic = 1
}
